# Update Function export header
#
# The "AUTO-TAGGING" sheet had a duplicated "Issue 8th" header in the merged
# Function-export banner (both P2 and Q2 read "Issue 8th"), which pushed every
# later "Issue N" header out of alignment. Remove the stray column so the
# header sequence (Issue 8th, 9th, 10th, Issue, Comment) lines up again, and
# drop the two leftover stray values that had spilled into column K on the
# detail rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AUTO-TAGGING")

# Delete the whole duplicate-header column; this shifts R:V left by one,
# fixing row 1 (the merged banner) and row 2 (the column headers) together.
$ws.Columns("Q").Delete()

# Remove the two stray "2" values that had leaked into column K.
$ws.Range("K14").Clear()
$ws.Range("K15").Clear()
